# CRM-4377: Add Inactive SF also in SF Document List Download with Active
# and Inactive status.
#
# The SF_List_Template worksheet gains two new columns:
#   - "Company Name" ({vendor:company_name}) inserted right after column A
#     (pushing every existing column one slot to the right), and
#   - "Status" ({vendor:active_status}) appended as the new last column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new "Company Name" column at B, shifting the rest right ---
$ws.Columns.Item(2).Insert()

# Fill the new column's data row first, then header (matches the order the
# workbook's shared-string table was actually populated in).
$ws.Range("B2").Value = "{vendor:company_name}"
$ws.Range("B2").Font.Bold = $false
$ws.Range("B2").HorizontalAlignment = -4108  # xlCenter

# --- Append the new "Status" column after the existing last column (AF) ---
$ws.Range("AG1").Value = "Status"
$ws.Range("AG1").Font.Bold = $true
$ws.Range("AG1").HorizontalAlignment = -4108  # xlCenter

$ws.Range("AG2").Value = "{vendor:active_status}"
$ws.Range("AG2").Font.Bold = $false
$ws.Range("AG2").HorizontalAlignment = -4108  # xlCenter

# Finally, set the "Company Name" header text.
$ws.Range("B1").Value = "Company Name"
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4108  # xlCenter

# --- Cosmetic column widths for the two brand-new columns ---
$ws.Columns.Item(32).ColumnWidth = 33.05
$ws.Columns.Item(33).ColumnWidth = 20.5

# --- Restore a plain selection (no frozen/scrolled top-left cell) ---
$ws.Range("B5").Select()
